# PNAD 2009 - correção nos dados e inicio da analise
# The sheet has two "section header" rows (row 5: "situação do domicílio" and,
# after the first deletion shifts things up, the original row 8:
# "grandes regiões e unidades da federação") that only carried a label in
# column A with no data. These two rows are removed entirely, which shifts
# every row below them up (cumulatively by 2), filling columns B:F of what
# become rows 5 and 8 with numeric data that previously lived two rows below.
# The last two rows of the old sheet (old rows 39/40) therefore drop off the
# bottom, shrinking the used range from F40 to F38.
#
# Also, the header label in B2 ("unnamed: 1_level_1") is corrected to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mislabeled header
$ws.Range("B2").Value = "total"

# Remove the "situação do domicílio" section-header-only row (row 5).
# Everything below shifts up by one row.
$ws.Rows("5").Delete()

# After the above delete, the former "grandes regiões e unidades da
# federação" section-header-only row (originally row 8) now sits at row 7.
# Remove it too; everything below shifts up by one more row.
$ws.Rows("7").Delete()
